$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data in A1:A11
$ws.Range("A1:A11").ClearContents()

# Set the new value
$ws.Range("A12").Value = "assas"

# Move the active selection to L12 (matching the final selection in the file)
$ws.Range("L12").Select()
